$d = $word.ActiveDocument

function Set-ParagraphXml($anchorText, $xml) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -like "*$anchorText*") {
            $rng = $p.Range
            $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + $xml + '</pkg:xmlData></pkg:part></pkg:package>'
            $rng.InsertXML($pkg)
            return $true
        }
    }
    return $false
}

$wDoc = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$wDocClose = '</w:document>'

# 1. Remove the stray _GoBack bookmark that sat after "Udsæt alarm, Kalibrér."
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. "Tiden der går ... til at grafen vises skal være 1 sek. med en tolerance på +/-15%."
#    -> split/reword, drop the tolerance clause.
$p8 = '<w:body><w:p w:rsidR="00266D10" w:rsidRDefault="00266D10" w:rsidP="00FC2CBB"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr></w:pPr><w:r w:rsidRPr="00DC2207"><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t>Tiden der går fra der er trykket</w:t></w:r><w:r w:rsidR="00DC2207"><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t xml:space="preserve"> på start-knappen</w:t></w:r><w:r w:rsidRPr="00DC2207"><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t xml:space="preserve"> til</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t xml:space="preserve"> at grafen vises må maksimalt</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t xml:space="preserve"> være 1 sek</w:t></w:r><w:r w:rsidRPr="00DC2207"><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body>'
Set-ParagraphXml "til at grafen vises" ($wDoc + $p8 + $wDocClose)

# 3. "GUI krav" -> "GUI " + "krav" wrapped with proofErr spell markers, bookmark moves here.
$p10 = '<w:body><w:p w:rsidR="00D2135A" w:rsidRPr="00D2135A" w:rsidRDefault="00D2135A" w:rsidP="00D2135A"><w:pPr><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t xml:space="preserve">GUI </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t>krav</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p></w:body>'
Set-ParagraphXml "GUI krav" ($wDoc + $p10 + $wDocClose)

# 4. "GUI'en skal se således ud" -> "GUI'e" + "rne" + " skal se således ud"
$p12 = '<w:body><w:p w:rsidR="00D2135A" w:rsidRDefault="00D2135A" w:rsidP="00DC2207"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t>GUI’e</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t>rne</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Segoe UI"/></w:rPr><w:t xml:space="preserve"> skal se således ud</w:t></w:r></w:p></w:body>'
Set-ParagraphXml "GUI’en skal se således ud" ($wDoc + $p12 + $wDocClose)

# 5. "Billedet indsættes til sidst" -> "Billederne indsættes her - skitser"
$d.Content.Find.Execute("Billedet indsættes til sidst", $false, $false, $false, $false, $false, $true, 1, $false, "Billederne indsættes her - skitser", 2) | Out-Null
